$wb = $excel.ActiveWorkbook

$wsFBS = $wb.Worksheets.Item("FBS")
$wsFBS.Range("N2").Value = 'NE'
$wsFBS.Range("O2").Value = 37.76
$wsFBS.Range("P2").Value = 12.5
$wsFBS.Range("S2").Value = -2
$wsFBS.Range("U2").Value = 2.2
$wsFBS.Range("AK2").Value = '2024-12-05T16:21:15.923721'
$wsFBS.Range("N3").Value = 'NNW'
$wsFBS.Range("P3").Value = 6.9
$wsFBS.Range("R3").Value = 1.2
$wsFBS.Range("U3").Value = -3.4
$wsFBS.Range("AK3").Value = '2024-12-05T16:21:15.923721'
$wsFBS.Range("M4").Value = 'ESE'
$wsFBS.Range("N4").Value = 'SSW'
$wsFBS.Range("O4").Value = 24.32
$wsFBS.Range("P4").Value = 1.3
$wsFBS.Range("Q4").Value = 'ESE'
$wsFBS.Range("S4").Value = -0.71
$wsFBS.Range("T4").Value = -0.96
$wsFBS.Range("U4").Value = -4.5
$wsFBS.Range("AB4").Value = 4.5
$wsFBS.Range("AF4").Value = -0.5
$wsFBS.Range("AK4").Value = '2024-12-05T16:21:15.923721'
$wsFBS.Range("A5").Value = 'Western Kentucky @ Jacksonville State'
$wsFBS.Range("D5").Value = 'Low'
$wsFBS.Range("E5").Value = 'E-W'
$wsFBS.Range("F5").Value = 'High'
$wsFBS.Range("G5").Value = 'N'
$wsFBS.Range("H5").Value = 43.63323969999999
$wsFBS.Range("I5").Value = 63.15
$wsFBS.Range("J5").Value = 59.04
$wsFBS.Range("K5").Value = 4.8
$wsFBS.Range("L5").Value = 1947
$wsFBS.Range("N5").Value = 'SSE'
$wsFBS.Range("O5").Value = 34.16
$wsFBS.Range("P5").Value = 5.1
$wsFBS.Range("Q5").Value = 'SSE'
$wsFBS.Range("S5").Value = 0
$wsFBS.Range("T5").Value = 0
$wsFBS.Range("U5").Value = 0.3
$wsFBS.Range("V5").Value = '33.8201052, -85.76647'
$wsFBS.Range("Z5").Value = -110
$wsFBS.Range("AA5").Value = -3.5
$wsFBS.Range("AF5").Value = 0.5
$wsFBS.Range("AK5").Value = '2024-12-05T16:21:15.923721'
$wsFBS.Range("A6").Value = 'UNLV @ Boise State'
$wsFBS.Range("D6").Value = 'High'
$wsFBS.Range("E6").Value = 'N-S'
$wsFBS.Range("F6").Value = 'Med'
$wsFBS.Range("G6").Value = 'E'
$wsFBS.Range("H6").Value = $null
$wsFBS.Range("I6").Value = 53.65
$wsFBS.Range("J6").Value = 70.04
$wsFBS.Range("K6").Value = 6.8
$wsFBS.Range("L6").Value = 1970
$wsFBS.Range("O6").Value = 30.02
$wsFBS.Range("P6").Value = 3
$wsFBS.Range("T6").Value = -0.25
$wsFBS.Range("U6").Value = -3.8
$wsFBS.Range("V6").Value = '43.6028839, -116.1958882'
$wsFBS.Range("Z6").Value = -105
$wsFBS.Range("AA6").Value = -4
$wsFBS.Range("AB6").Value = -4
$wsFBS.Range("AF6").Value = 0
$wsFBS.Range("AK6").Value = '2024-12-05T16:21:15.923721'
$wsFBS.Range("M7").Value = 'W'
$wsFBS.Range("N7").Value = 'W'
$wsFBS.Range("O7").Value = 50.66
$wsFBS.Range("P7").Value = 7.7
$wsFBS.Range("Q7").Value = 'W'
$wsFBS.Range("R7").Value = 0.4
$wsFBS.Range("U7").Value = -2.2
$wsFBS.Range("Y7").Value = 57.5
$wsFBS.Range("Z7").Value = 100
$wsFBS.Range("AE7").Value = -0.0170940170940171
$wsFBS.Range("AK7").Value = '2024-12-05T16:21:15.923721'

$wsOther = $wb.Worksheets.Item("Other")
$wsOther.Range("O2").Value = 'ENE'
$wsOther.Range("P2").Value = 'ENE'
$wsOther.Range("Q2").Value = 44.54
$wsOther.Range("R2").Value = 11.3
$wsOther.Range("S2").Value = 'ENE'
$wsOther.Range("P3").Value = 'ENE'
$wsOther.Range("Q3").Value = 53.6
$wsOther.Range("R3").Value = 5.2
$wsOther.Range("O4").Value = 'SW'
$wsOther.Range("P4").Value = 'WSW'
$wsOther.Range("Q4").Value = 52.16
$wsOther.Range("R4").Value = 7.8
$wsOther.Range("S4").Value = 'WSW'
$wsOther.Range("T4").Value = 0
$wsOther.Range("Q5").Value = 60.14000000000001
$wsOther.Range("R5").Value = 4.5
